# Auto-generated edit script applying cryptos.xlsx price/volume/coin updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.198.00"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "2.951.97"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'378.49"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").Value = "'102.23"
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("E7").Value = "  -0.52%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("D10").Value = "'36.51"
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "'0.0840"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").Value = "3.417.22"
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").Value = "'18.00"
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "2.944.63"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").Value = "'0.987"
$ws.Range("E17").Value = "  +5.47%  "
$ws.Range("D18").Value = "51.114.49"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("E19").Value = "  -5.20%  "
$ws.Range("D20").Value = "'7.20"
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("D21").Value = "'12.60"
$ws.Range("E21").Value = "  -2.55%  "
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("D23").Value = "'68.39"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").Value = "'261.73"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("E25").Value = "  +2.73%  "
$ws.Range("D26").Value = "'8.39"
$ws.Range("E26").Value = "  +13.38%  "
$ws.Range("E27").Value = "  +8.08%  "
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").Value = "'4.11"
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.113"
$ws.Range("E30").Value = "  +11.81%  "
$ws.Range("B31").Value = "Dai"
$ws.Range("C31").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "'33.90"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'50.48"
$ws.Range("E35").Value = "  -2.45%  "
$ws.Range("D36").Value = "'2.05"
$ws.Range("E36").Value = "  -2.18%  "
$ws.Range("D37").Value = "'0.0450"
$ws.Range("E37").Value = "  +5.69%  "
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("E39").Value = "  -1.40%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").Value = "'2.57"
$ws.Range("E41").Value = "  -0.80%  "
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("D44").Value = "'121.86"
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("D45").Value = "'21.01"
$ws.Range("E45").Value = "  -3.74%  "
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").Value = "'0.273"
$ws.Range("E47").Value = "  +1.09%  "
$ws.Range("E48").Value = "  +2.13%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'3.22"
$ws.Range("E49").Value = "  +1.81%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.003.27"
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("D51").Value = "'0.0335"
$ws.Range("E51").Value = "  +5.50%  "
